$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new TLC59711 LED driver row beneath the existing example row.
# Cells are populated in the same order the shared strings table ends up
# recording them in, so new <si> entries land in the expected sequence.
$ws.Range("G10").Value = "595-TLC59711PWP"
$ws.Range("H10").Value = "TI"
$ws.Range("I10").Value = "TLC59711PWP"
$ws.Range("J10").Value = "4.12/3.37/2.47"
$ws.Range("E10").Value = "HTSSOP20"
$ws.Range("D10").Value = "U102"
$ws.Range("C10").Value = "TLC59711"
$ws.Range("B10").Value = "6 min."
$ws.Range("A10").Value = "LED driver"
$ws.Range("F10").Value = "Mouser"

# Update the selected cell, matching the saved view state
[void]$ws.Range("D11").Select()
